$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace the product catalog contents (Sheet1, A1:D3) ---
# Row 1 (was the "Name/Description/Price/Image" header row; now data row for EndStrips)
$ws.Range("A1").Value = "EndStrips"
$ws.Range("B1").Value = "Revolutionary micro-dosing hypoglycemic system. Developed in house by EndT1."
$ws.Range("C1").Value = 10
$ws.Range("D1").Value = "img/endstrips.png"

# Row 2 (EndCaps product)
$ws.Range("A2").Value = "EndCaps"
$ws.Range("B2").Value = "Extended release glucose capsules. Never worry about hypoglycemia at night again!"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "img/endcaps.png"

# Row 3 (test row)
$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "asdfasd"
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = "img/endcaps.png"

# --- Formatting ---
# The old sheet bolded the header row (A1:D1) via font style; the new layout has
# no header row (it's plain data now), so drop the bold.
$ws.Range("A1:D1").Font.Bold = $false

# The "price" column (C1:C2, value 10) now carries a 2-decimal number format.
$ws.Range("C1:C2").NumberFormat = "0.00"

# --- Column widths ---
# Target widths (per-character units) are 8.375 / 71 / 5.375 / 26.125; this
# runtime's ColumnWidth setter adds a fixed 5/6-character pad before storing,
# so we back that off to land as close as representable to the target widths.
$ws.Columns.Item(1).ColumnWidth = 7.5
$ws.Columns.Item(2).ColumnWidth = 70.16666666666667
$ws.Columns.Item(3).ColumnWidth = 4.5
$ws.Columns.Item(4).ColumnWidth = 25.333333333333336

# --- Selection moves from B8 to D3 ---
[void]$ws.Range("D3").Select()

Write-Output "done"
